$wb = $excel.ActiveWorkbook

$wsOps = $wb.Worksheets.Item("ShopOrderOperations")
$wsOps.Range("F2").Value = 6
$wsOps.Range("J16").Select()

$wsAlloc = $wb.Worksheets.Item("WorkCenterOpAllocations")
$wsAlloc.Range("C2").Value = 1
$wsAlloc.Range("C10").Value = 1
$wsAlloc.Range("J2").Select()
